$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45961, 6941.95924550394, 5769.21464309394, 3860, 2580.31782175543, 187.063852702057),
    @(3, 45962, 2042.57134452881, 3377.28986618857, 9004, 3144.22157429192, -103.437023313313),
    @(4, 45963, 2008.80229450428, 3493.28808195605, 9004, 3104.30517170886, -100.266947763962),
    @(5, 45964, 7982.44608295771, 6800.36500999569, 9004, 4211.71624739826, 83.6700523914146),
    @(6, 45965, 8215.37031519319, 7542.37676062885, 9004, 4384.54078630116, 121.788231122084),
    @(7, 45966, 8517.17560294576, 8008.00364360812, 9004, 4586.81288627443, 149.617355411773),
    @(8, 45967, 8525.26551938708, 8230.44233199449, 9004, 4588.76921430204, 158.967147762355),
    @(9, 45968, 8530.1901118937, 7461.3666833131, 9004, 4619.39747280756, 128.198506505028),
    @(10, 45969, 3118.30130226871, 4920.61223718882, 9004, 4308.06435720446, 9.36152476638684),
    @(11, 45970, 2984.26537892464, 4793.97691111794, 9004, 4300.79824928305, 3.78229835004125),
    @(12, 45971, 9278.26748203409, 8023.91017910924, 9004, 5247.5641376462, 177.81142986481),
    @(13, 45972, 9278.26878400516, 8296.80888643121, 9004, 5247.56790175123, 189.182366174268),
    @(14, 45973, 9278.26878400516, 8380.82245606453, 9004, 5247.56790175123, 192.682931575657),
    @(15, 45974, 9278.26878400516, 8342.57698159789, 9004, 5247.56790175123, 191.089370139547)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
